$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = [double]"22.51000000000008"
$ws.Range("H2").Value = [double]"1.787079315292003e-16"
$ws.Range("K2").Value = [double]"46.5280221252621"
$ws.Range("L2").Value = "[43.70992394338282, 49.346120307141376]"
$ws.Range("O2").Value = [double]"1.616395018964117"
$ws.Range("P2").Value = "[1.5535002711445012, 1.6792897667837332]"
$ws.Range("S2").Value = [double]"51.3326648698968"
$ws.Range("T2").Value = "[49.41068548570021, 53.25464425409339]"
$ws.Range("W2").Value = [double]"16.7191391391392"
$ws.Range("X2").Value = [double]"16.49381381381387"
$ws.Range("Y2").Value = [double]"16.94446446446453"

# Row 3
$ws.Range("E3").Value = [double]"22.40000000000006"
$ws.Range("H3").Value = [double]"1.787079315292003e-16"
$ws.Range("I3").Value = [double]"0.6137674231351554"
$ws.Range("K3").Value = [double]"45.86270272222085"
$ws.Range("L3").Value = "[42.478903897452724, 49.24650154698898]"
$ws.Range("O3").Value = [double]"1.113237036407194"
$ws.Range("P3").Value = "[1.0377633390236563, 1.1887107337907326]"
$ws.Range("S3").Value = [double]"52.88735303215507"
$ws.Range("T3").Value = "[51.06867226326236, 54.70603380104778]"
$ws.Range("W3").Value = [double]"18.43123123123128"
$ws.Range("X3").Value = [double]"18.16216216216221"
$ws.Range("Y3").Value = [double]"18.70030030030035"
